$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Collapse the "Optimized and accelerated model training..." bullet
#    (originally split across many runs) into the new single-run text.
#    Find/Replace naturally merges the matched runs into one run that
#    carries the formatting of the first matched run, which is exactly
#    what the target OOXML shows.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Optimized and accelerated model training by 30% by tuning learning rate and optimizer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Accelerated network training by 30% training model parallelly with JAX", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Append a *new, separate* run containing ", Kubernetes" right after
#    "PyTorch, TensorFlow, Keras, Hadoop, AWS" without merging it into
#    the existing run (InsertAfter would merge same-formatted adjacent
#    text into the existing run). We build the text in a throwaway
#    paragraph at the end of the story, cut it back out (so the run
#    boundary survives), tidy up the scratch paragraph, then paste the
#    cut fragment at the target location -- Paste preserves the run as
#    its own element instead of folding it into neighboring text.
# ------------------------------------------------------------------

# 2a. Create a scratch paragraph at the very end of the document body.
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertParagraphAfter()

# 2b. Type the new fragment into that scratch paragraph.
$scratch = $d.Content
$scratch.Collapse(0)
$scratch.InsertAfter(", Kubernetes")

# 2c. Select it and Cut it to the clipboard (removes the text, leaves
#     the now-empty scratch paragraph behind).
$scratch.Select()
$word.Selection.Cut() | Out-Null

# 2d. Remove the now-empty scratch paragraph (delete its preceding
#     paragraph mark together with the (empty) paragraph range).
$lastPara = $d.Paragraphs.Last
$scratchRange = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$scratchRange.Delete() | Out-Null

# 2e. Find the insertion point right after "...Hadoop, AWS" and paste
#     the cut fragment there as its own run.
$target = $d.Content
$target.Find.Execute(
    "PyTorch, TensorFlow, Keras, Hadoop, AWS",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Collapse(0)
$target.Select()
$word.Selection.Paste() | Out-Null
